# Updates cryptocurrency price/volume figures in columns D (Price) and E (Volume(1h))
# for the rows whose source data changed in this refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.554.90"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "1.640.87"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'308.33"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").Value = "'52.69"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("D9").Value = "'0.3668"
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("D10").Value = "'1.279"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").Value = "'0.08198"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "'0.9997"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "'23.02"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").Value = "'6.669"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").Value = "'7.434"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "1.642.55"
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("D18").Value = "'95.00"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").Value = "'0.06925"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").Value = "'18.31"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").Value = "'6.578"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("D22").Value = "'0.9987"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "23.548.77"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").Value = "'12.88"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").Value = "'3.082"
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("D26").Value = "'2.420"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "'151.51"
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("D29").Value = "'5.362"
$ws.Range("E29").Value = "  +2.58%  "
$ws.Range("D30").Value = "'136.00"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("D31").Value = "'2.390"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").Value = "1.827.62"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("D33").Value = "'6.837"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "'0.9784"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").Value = "'0.02842"
$ws.Range("E35").Value = "  +4.40%  "
$ws.Range("D36").Value = "'10.40"
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("D37").Value = "'0.07385"
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("D38").Value = "'0.2556"
$ws.Range("D39").Value = "'6.216"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("D40").Value = "'0.08898"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").Value = "'1.387"
$ws.Range("E41").Value = "  +1.82%  "
$ws.Range("D42").Value = "'0.7128"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "'16.41"
$ws.Range("E43").Value = "  +6.04%  "
$ws.Range("D44").Value = "'12.56"
$ws.Range("E44").Value = "  +1.19%  "
$ws.Range("D45").Value = "'0.6569"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'2.351"
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("D48").Value = "'0.9987"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "'130.21"
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").Value = "'1.210"
$ws.Range("E51").Value = "  +0.35%  "
